$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 2251
$wsExhibit.Range("F5").Value = 1713
$wsExhibit.Range("F8").Value = 781

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2251
$wsAll.Range("F5").Value = 1713
$wsAll.Range("F9").Value = 781
